$d = $word.ActiveDocument

$replacements = @(
    @("21÷7=", "26÷9="),
    @("66÷8=", "35÷6="),
    @("86÷9=", "60÷8="),
    @("28÷6=", "88÷5="),
    @("26÷7=", "51÷7="),
    @("29÷4=", "30÷8="),
    @("48÷3=", "14÷8="),
    @("36÷7=", "74÷7="),
    @("16÷4=", "50÷6="),
    @("44÷4=", "24÷5="),
    @("11÷6=", "18÷3="),
    @("66÷2=", "64÷9="),
    @("41÷3=", "92÷5="),
    @("27÷5=", "76÷2="),
    @("45÷2=", "31÷5="),
    @("89÷7=", "10÷6="),
    @("16÷5=", "10÷4="),
    @("41÷5=", "52÷4="),
    @("58÷8=", "96÷4="),
    @("52÷7=", "57÷7="),
    @("76÷6=", "69÷3="),
    @("72÷7=", "25÷6="),
    @("53÷9=", "42÷6="),
    @("93÷8=", "62÷2="),
    @("88÷3=", "88÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
